# Applies the "Home Page.docx" revision:
#   1. Bold the "My Portfolio" heading paragraph.
#   2. Insert a new paragraph (about long-term career goals) between the
#      existing blank paragraph and the "About Me" heading, followed by a
#      new blank paragraph.
#   3. Bold the "About Me" heading paragraph.

$d = $word.ActiveDocument

function Set-ParagraphBold($para) {
    $para.Range.Bold = 1
    $para.Range.Font.BoldBi = 1
}

# --- 1. Bold "My Portfolio" -------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
Set-ParagraphBold $titlePara

# --- 2. Insert the new "Long term..." paragraph -----------------------------
# The paragraph is built out of many separate runs (matching the source
# revision) by writing each text segment into its own temporary paragraph and
# then deleting the intervening paragraph marks to splice them back together
# into one paragraph without merging the runs.

$segments = @(
    'Long term, I would like to become specialized in some type of field',
    ' with a Graduate degree',
    ', right now leaning towards',
    ' something related to',
    ' PCB design, FPG',
    'As',
    ', or possibly something related to digital signals processing or radio frequency devices. ',
    'I hope to gain a better understanding in each of these fields as I complete my last year of ',
    'my Undergraduate degree ',
    'and my Masters at Iowa State University. After graduation, I am ',
    'planning to stay within the Midwest so that I can be closer to family. ',
    'I am looking forward to entering industry to apply the skills I have learned in my classes, projects, and internships towards something that can help improve the lives of people daily. '
)

# Anchor: the blank paragraph that currently sits right before "About Me".
$anchor = $d.Paragraphs.Item(3)
$firstNewIndex = $anchor.Index + 1

$lastPara = $anchor
foreach ($seg in $segments) {
    $lastPara.Range.InsertParagraphAfter()
    $lastPara = $d.Paragraphs.Item($lastPara.Index + 1)
    $lastPara.Range.InsertAfter($seg)
}

# Splice the per-segment paragraphs into a single paragraph, one merge per
# extra segment, keeping each segment as its own run.
for ($k = 1; $k -lt $segments.Count; $k++) {
    $p = $d.Paragraphs.Item($firstNewIndex)
    $endOfP = $p.Range.End
    $markRange = $d.Range($endOfP - 1, $endOfP)
    $markRange.Delete()
}

# --- New blank paragraph between the new text and "About Me" ---------------
# InsertParagraphAfter() always leaves a placeholder empty run behind; type a
# throwaway character and delete it again so the paragraph collapses back to
# a truly empty <w:p/>, matching the blank paragraphs already in the doc.
$newTextPara = $d.Paragraphs.Item($firstNewIndex)
$newTextPara.Range.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Item($firstNewIndex + 1)
$blankPara.Range.InsertAfter("X")
$blankPara2 = $d.Paragraphs.Item($firstNewIndex + 1)
$placeholderRange = $d.Range($blankPara2.Range.Start, $blankPara2.Range.Start + 1)
$placeholderRange.Delete()

# --- 3. Bold "About Me" ------------------------------------------------------
$aboutIndex = $firstNewIndex + 2
$aboutPara = $d.Paragraphs.Item($aboutIndex)
Set-ParagraphBold $aboutPara
